$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")
$wsOut = $wb.Worksheets.Item("ProductLoanOutput")

$ws.Range("B1").Value = "2620-MS-EI-DB-DL-REC-NOCOM-RNI-CTPD-DL-MD-TR-2-DATE-VAR-INST-OVERDUE-FEE-%INT-LESS-AMT-1st"
$ws.Range("B2").Value = "262d"
$wsOut.Range("B1").Value = "2620-MS-EI-DB-DL-REC-NOCOM-RNI-CTPD-DL-MD-TR-2-DATE-VAR-INST-OVERDUE-FEE-%INT-LESS-AMT-1st"

$ws.Range("B3").Select()
